$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2023-08-26
$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 7.367362464249472

# Row 3: 2023-08-23
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 4.371470058157054

# Row 4: 2023-08-12
$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 4.371470058157054

# Row 5: 2023-07-04
$ws.Range("B5").Value = 0.006876353814593728
$ws.Range("C5").Value = 0.05231270169004087
$ws.Range("D5").Value = 2938.103010863317
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 3185.147450612923
